$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old date -> new date (shift the week forward by one, 2023-09-18..23 -> 2023-09-25..30)
$dateMap = @{
    "2023-09-18" = "2023-09-25"
    "2023-09-19" = "2023-09-26"
    "2023-09-20" = "2023-09-27"
    "2023-09-21" = "2023-09-28"
    "2023-09-22" = "2023-09-29"
    "2023-09-23" = "2023-09-30"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Text
    if ($dateMap.ContainsKey($val)) {
        $cell.Value = $dateMap[$val]
    }
}

# Update the active selection shown in the sheet view
$ws.Range("B43").Select()
